# Update odds values for the week's games (FlashScore 2024-10-07 dataset)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Goias vs Santos
$ws.Range("Z2").Value = 29
$ws.Range("AG2").Value = 6.5
$ws.Range("AW2").Value = 4.5
$ws.Range("AX2").Value = 17

# Row 5 - Fortaleza vs Jaguares de Cordoba
$ws.Range("G5").Value = 1.55
$ws.Range("U5").Value = 2.38
$ws.Range("V5").Value = 1.53
$ws.Range("Z5").Value = 10
$ws.Range("AI5").Value = 23
$ws.Range("AO5").Value = 8
$ws.Range("AQ5").Value = 26

# Row 8 - Macara vs Delfin
$ws.Range("G8").Value = 1.95
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 2.75
$ws.Range("L8").Value = 4.5
$ws.Range("Z8").Value = 17
$ws.Range("AE8").Value = 15
$ws.Range("AN8").Value = 4
$ws.Range("AW8").Value = 5.5
$ws.Range("BB8").Value = 251
